$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name value (added a hyphen after "246")
$wsInput.Range("B1").Value = "246-MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"
$wsOutput.Range("B1").Value = "246-MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"

# Select B1 on the input sheet, then activate output sheet and select B1 there too,
# leaving ProductLoanOutput as the active/visible tab.
$wsInput.Activate()
$wsInput.Range("B1").Select() | Out-Null

$wsOutput.Activate()
$wsOutput.Range("B1").Select() | Out-Null
